{"js": "// Update the date line and every two-digit \u00d7 two-digit multiplication\n// problem in the table to the new set of values, matching the target\n// OOXML diff exactly (plain text substitutions, formatting untouched).\nconst replacements = [\n  [\"2024-05-30 Thursday\", \"2024-05-31 Friday\"],\n  [\"38\u00d798=\", \"99\u00d762=\"],\n  [\"79\u00d713=\", \"43\u00d722=\"],\n  [\"69\u00d774=\", \"76\u00d776=\"],\n  [\"19\u00d798=\", \"83\u00d738=\"],\n  [\"87\u00d778=\", \"66\u00d786=\"],\n  [\"62\u00d770=\", \"68\u00d720=\"],\n  [\"37\u00d758=\", \"54\u00d798=\"],\n  [\"23\u00d748=\", \"19\u00d790=\"],\n  [\"60\u00d742=\", \"80\u00d769=\"],\n  [\"21\u00d766=\", \"45\u00d734=\"],\n  [\"49\u00d747=\", \"36\u00d718=\"],\n  [\"56\u00d729=\", \"46\u00d763=\"],\n  [\"50\u00d766=\", \"22\u00d786=\"],\n  [\"72\u00d741=\", \"18\u00d783=\"],\n  [\"46\u00d755=\", \"76\u00d742=\"],\n  [\"28\u00d776=\", \"69\u00d716=\"],\n  [\"14\u00d762=\", \"42\u00d774=\"],\n  [\"63\u00d753=\", \"34\u00d783=\"],\n  [\"96\u00d779=\", \"32\u00d797=\"],\n  [\"15\u00d771=\", \"92\u00d721=\"],\n  [\"55\u00d759=\", \"18\u00d780=\"],\n  [\"27\u00d768=\", \"31\u00d796=\"],\n  [\"74\u00d759=\", \"78\u00d775=\"],\n  [\"34\u00d744=\", \"30\u00d798=\"],\n  [\"33\u00d779=\", \"80\u00d770=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text to replace: \"${oldText}\"`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and every two-digit x two-digit multiplication\n# problem in the table to the new set of values, matching the target\n# OOXML diff exactly (plain text substitutions, formatting untouched).\n$replacements = @(\n    @{ Old = \"2024-05-30 Thursday\"; New = \"2024-05-31 Friday\" },\n    @{ Old = \"38\u00d798=\"; New = \"99\u00d762=\" },\n    @{ Old = \"79\u00d713=\"; New = \"43\u00d722=\" },\n    @{ Old = \"69\u00d774=\"; New = \"76\u00d776=\" },\n    @{ Old = \"19\u00d798=\"; New = \"83\u00d738=\" },\n    @{ Old = \"87\u00d778=\"; New = \"66\u00d786=\" },\n    @{ Old = \"62\u00d770=\"; New = \"68\u00d720=\" },\n    @{ Old = \"37\u00d758=\"; New = \"54\u00d798=\" },\n    @{ Old = \"23\u00d748=\"; New = \"19\u00d790=\" },\n    @{ Old = \"60\u00d742=\"; New = \"80\u00d769=\" },\n    @{ Old = \"21\u00d766=\"; New = \"45\u00d734=\" },\n    @{ Old = \"49\u00d747=\"; New = \"36\u00d718=\" },\n    @{ Old = \"56\u00d729=\"; New = \"46\u00d763=\" },\n    @{ Old = \"50\u00d766=\"; New = \"22\u00d786=\" },\n    @{ Old = \"72\u00d741=\"; New = \"18\u00d783=\" },\n    @{ Old = \"46\u00d755=\"; New = \"76\u00d742=\" },\n    @{ Old = \"28\u00d776=\"; New = \"69\u00d716=\" },\n    @{ Old = \"14\u00d762=\"; New = \"42\u00d774=\" },\n    @{ Old = \"63\u00d753=\"; New = \"34\u00d783=\" },\n    @{ Old = \"96\u00d779=\"; New = \"32\u00d797=\" },\n    @{ Old = \"15\u00d771=\"; New = \"92\u00d721=\" },\n    @{ Old = \"55\u00d759=\"; New = \"18\u00d780=\" },\n    @{ Old = \"27\u00d768=\"; New = \"31\u00d796=\" },\n    @{ Old = \"74\u00d759=\"; New = \"78\u00d775=\" },\n    @{ Old = \"34\u00d744=\"; New = \"30\u00d798=\" },\n    @{ Old = \"33\u00d779=\"; New = \"80\u00d770=\" },\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $pair.New\n\n    $found = $find.Execute(\n        $pair.Old,   # FindText\n        $true,       # MatchCase\n        $false,      # MatchWholeWord\n        $false,      # MatchWildcards\n        $false,      # MatchSoundsLike\n        $false,      # MatchAllWordForms\n        $true,       # Forward\n        1,           # Wrap (wdFindContinue)\n        $false,      # Format\n        $pair.New,   # ReplaceWith\n        2            # Replace (wdReplaceAll)\n    )\n\n    if (-not $found) {\n        throw \"Could not find text to replace: $($pair.Old)\"\n    }\n}\n"}
